$d = $word.ActiveDocument

# 1. Update the date from 2022 to 2023
$d.Content.Find.Execute("October 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "October 2023", 2)

# 2. Update the intro paragraph: Tutorial Question 3 -> Tutorial Question 4
#    (leave "Problem Sheet 3" unchanged)
$d.Content.Find.Execute("Tutorial Question 3 off of Problem Sheet 3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Tutorial Question 4 off of Problem Sheet 3", 2)

# 3. Update the example heading: PS3 Question 3 -> PS3 Question 4
$d.Content.Find.Execute("Example 1 (PS3 Question 3)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Example 1 (PS3 Question 4)", 2)

# 4. Update the Q&A text: 3c) -> 4c)
$d.Content.Find.Execute("Why is 3c) done in this way?", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Why is 4c) done in this way?", 2)
